$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> Masum Electronics / DSR-0248 retailer record
$ws.Cells.Item(2,2).Value  = "DSR-0248"
$ws.Cells.Item(2,3).Value  = "Masum Electronics"
$ws.Cells.Item(2,4).Value  = "Khidilpur"
$ws.Cells.Item(2,5).Value  = "Md Masum Ali"
$ws.Cells.Item(2,9).Value  = "Md Masum Ali"
$ws.Cells.Item(2,10).Value = 1632005795
$ws.Cells.Item(2,11).Value = "Natore"
$ws.Cells.Item(2,12).Value = "Baraigram"
$ws.Cells.Item(2,14).Value = "Khidilpur Bazar, Mokhura,Baraigram"
$ws.Cells.Item(2,16).Value = 1632005795
$ws.Cells.Item(2,20).Value = 1632005795

# Row 3 -> SR Telecom / DSR-0351 retailer record
$ws.Cells.Item(3,2).Value  = "DSR-0351"
$ws.Cells.Item(3,3).Value  = "SR Telecom"
$ws.Cells.Item(3,4).Value  = "Jholmolia"
$ws.Cells.Item(3,5).Value  = "Maqsudur Rahman"
$ws.Cells.Item(3,9).Value  = "Maqsudur Rahman"
$ws.Cells.Item(3,10).Value = 1768927219
$ws.Cells.Item(3,11).Value = "Rajshahi"
$ws.Cells.Item(3,12).Value = "Puthia"
$ws.Cells.Item(3,14).Value = "Jholmolia , Puthia, Rajshahi"
$ws.Cells.Item(3,16).Value = 1768927219
$ws.Cells.Item(3,20).Value = 1768927219

# Update the selected / top-left view state to match the saved workbook view
$sel = $ws.Range("E14").Select()
